$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 106 ---
$ws.Range('E106').Value = 'Attempt predict3dunet with same patch & stride as during training.'
$ws.Range('F106').Value = 'Success (no error)'
$ws.Range('G106').Value = 'Fail. Patch shape invalid error.'
$ws.Range('H106').Value = 0
$ws.Range('I106').Value = 1
$ws.Range('J106').Value = 'error, invalid patch shape &/ invalid stride shape'
$ws.Range('K106').Value = 0
$ws.Range('L106').Value = 'NA'

# --- Row 107 ---
$ws.Range('A107').Value = '231225-0'
$ws.Range('B107').Value = 'predict3dunet 1.6.0'
$ws.Range('C107').Value = 'ResidualUNet3D'
$ws.Range('D107').Value = '3DUnet_lightsheet_boundary'
$ws.Range('E107').Value = 'Attempt predict3dunet with same patch as during training, but stride shape such that there is a halo of 32 voxels in all dimensions.'
$ws.Range('F107').Value = 'Success (no error)'
$ws.Range('AM107').Value = 125
$ws.Range('AN107').Value = 1169
$ws.Range('AO107').Value = 414
$ws.Range('AP107').Value = 96
$ws.Range('AQ107').Value = 960
$ws.Range('AR107').Value = 256
$ws.Range('AT107').Value = 64
$ws.Range('AU107').Value = 928
$ws.Range('AV107').Value = 224
$ws.Range('BA107').Value = 'patch = same as used for training this model.'
$ws.Range('BB107').Value = 'stride = such that the halo is 32 in all dimensions.'

# --- Row 108 ---
$ws.Range('A108').Value = '231225-1'
$ws.Range('B108').Value = 'train3dunet 1.8.2'
$ws.Range('C108').Value = 'ResidualUNet3D'
$ws.Range('D108').Value = '3DUnet_lightsheet_boundary'
$ws.Range('E108').Value = 'Try out the new pytorch-3dunet 1.8.2 which has the ability of taking an arbitrary patch shape for ResidualUNet3D models. Maybe now predict3dunet works. Same patch as in most recent attempt.'
$ws.Range('F108').Value = 'TBD'
$ws.Range('G108').Value = 'TBD'
$ws.Range('H108').Value = 'TBD'
$ws.Range('I108').Value = 'TBD'
$ws.Range('J108').Value = 'TBD'
$ws.Range('K108').Value = 'TBD'
$ws.Range('L108').Value = 'TBD'
$ws.Range('M108').Value = 'TBD'
$ws.Range('N108').Value = 'TBD'
$ws.Range('O108').Value = 'TBD'
$ws.Range('P108').Value = 'TBD'
$ws.Range('Q108').Value = 'TBD'
$ws.Range('R108').Value = 'TBD'
$ws.Range('S108').Value = 'TBD'
$ws.Range('T108').Value = 'TBD'
$ws.Range('U108').Value = 'dataset07.0'
$ws.Range('V108').Value = 'TBD'
$ws.Range('W108').Value = 'TBD'
$ws.Range('X108').Value = 'TBD'
$ws.Range('Y108').Value = 'TBD'
$ws.Range('Z108').Value = 'TBD'
$ws.Range('AA108').Value = 'TBD'
$ws.Range('AB108').Value = 'TBD'
$ws.Range('AC108').Value = 'TBD'
$ws.Range('AD108').Value = 'TBD'
$ws.Range('AE108').Value = 'TBD'
$ws.Range('AF108').Value = 'TBD'
$ws.Range('AG108').Value = 'TBD'
$ws.Range('AH108').Value = 'TBD'
$ws.Range('AI108').Value = 'TBD'
$ws.Range('AJ108').Value = 'TBD'
$ws.Range('AK108').Value = 'TBD: formula is TBD'
$ws.Range('AL108').Value = 'NVIDIA A100-SXM4-80GB'
$ws.Range('AM108').Value = 125
$ws.Range('AN108').Value = 1169
$ws.Range('AO108').Value = 414
$ws.Range('AP108').Value = 96
$ws.Range('AQ108').Value = 960
$ws.Range('AR108').Value = 256
$ws.Range('AS108').Value = 'yes'
$ws.Range('AT108').Formula = '= FLOOR.MATH((AM108 - AP108) / 2)'
$ws.Range('AU108').Formula = '= FLOOR.MATH((AN108 - AQ108) / 2)'
$ws.Range('AV108').Formula = '= FLOOR.MATH((AO108 - AR108) / 2)'
$ws.Range('AW108').Value = 'no'
$ws.Range('AX108').Formula = '=AP108-AT108'
$ws.Range('AY108').Formula = '=AQ108-AU108'
$ws.Range('AZ108').Formula = '=AR108-AV108'
$ws.Range('BA108').Value = 'patch = some number = sum(2^i), with i>=4'
$ws.Range('BB108').Value = 'stride = floor (resolution - patch) / 2; for validation: stride = patch (like in Wolny''s config.yml)'
$ws.Range('BC108').Value = 'TBD'
$ws.Range('BD108').Value = 'TBD'
$ws.Range('BE108').Value = 'TBD'

# --- Row 109 ---
$ws.Range('A109').Value = 'TBD'
$ws.Range('B109').Value = 'TBD'
$ws.Range('C109').Value = 'TBD'
$ws.Range('D109').Value = 'TBD'
$ws.Range('E109').Value = 'TBD'
$ws.Range('F109').Value = 'TBD'
$ws.Range('G109').Value = 'TBD'
$ws.Range('H109').Value = 'TBD'
$ws.Range('I109').Value = 'TBD'
$ws.Range('J109').Value = 'TBD'
$ws.Range('K109').Value = 'TBD'
$ws.Range('L109').Value = 'TBD'
$ws.Range('M109').Value = 'TBD'
$ws.Range('N109').Value = 'TBD'
$ws.Range('O109').Value = 'TBD'
$ws.Range('P109').Value = 'TBD'
$ws.Range('Q109').Value = 'TBD'
$ws.Range('R109').Value = 'TBD'
$ws.Range('S109').Value = 'TBD'
$ws.Range('T109').Value = 'TBD'
$ws.Range('U109').Value = 'TBD'
$ws.Range('V109').Value = 'TBD'
$ws.Range('W109').Value = 'TBD'
$ws.Range('X109').Value = 'TBD'
$ws.Range('Y109').Value = 'TBD'
$ws.Range('Z109').Value = 'TBD'
$ws.Range('AA109').Value = 'TBD'
$ws.Range('AB109').Value = 'TBD'
$ws.Range('AC109').Value = 'TBD'
$ws.Range('AD109').Value = 'TBD'
$ws.Range('AE109').Value = 'TBD'
$ws.Range('AF109').Value = 'TBD'
$ws.Range('AG109').Value = 'TBD'
$ws.Range('AH109').Value = 'TBD'
$ws.Range('AI109').Value = 'TBD'
$ws.Range('AJ109').Value = 'TBD'
$ws.Range('AK109').Value = 'TBD: formula is TBD'
$ws.Range('AL109').Value = 'NVIDIA A100-SXM4-80GB'
$ws.Range('AM109').Value = 'TBD'
$ws.Range('AN109').Value = 'TBD'
$ws.Range('AO109').Value = 'TBD'
$ws.Range('AP109').Value = 'TBD'
$ws.Range('AQ109').Value = 'TBD'
$ws.Range('AR109').Value = 'TBD'
$ws.Range('AS109').Value = 'NA'
$ws.Range('AT109').Value = 'TBD'
$ws.Range('AU109').Value = 'TBD'
$ws.Range('AV109').Value = 'TBD'
$ws.Range('AW109').Value = 'NA'
$ws.Range('AX109').Formula = '=AP109-AT109'
$ws.Range('AY109').Formula = '=AQ109-AU109'
$ws.Range('AZ109').Formula = '=AR109-AV109'
$ws.Range('BA109').Value = 'TBD'
$ws.Range('BB109').Value = 'TBD'
$ws.Range('BC109').Value = 'TBD'
$ws.Range('BD109').Value = 'TBD'
$ws.Range('BE109').Value = 'TBD'

# --- Formatting fixups (row 108 mirrors row 106 formatting for these cells) ---
$ws.Range('AS106').Copy()
$ws.Range('AS108').PasteSpecial(-4122)
$ws.Range('AT106').Copy()
$ws.Range('AT108').PasteSpecial(-4122)
$excel.CutCopyMode = 0
